# The editorial note tracked in row 63 (sequence 62, "textAlign applied to
# span") has been addressed, so its Status moves from Open ("O") to
# Closed ("C").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E63").Value = "C"

# The sheet's AutoFilter previously filtered on both Status ("O") and
# Assigned to ("nigel"); the "Assigned to" criterion (column H, field 8)
# is removed, leaving only the Status filter. Re-applying AutoFilter on
# that field with no criteria clears it, and the engine recomputes which
# rows are hidden for us.
$af = $ws.AutoFilter
[void]$af.Range.AutoFilter(8)

# Leave the selection where the author's cursor ended up after the edit.
[void]$ws.Range("E107").Select()
